# citeis-auto TestData.xlsx : "added execute column in testdata.xlsx"
#
# Adds a new EXECUTE column (E) to the AppFlow sheet:
#   E1 = "EXECUTE" (same header look/style as A1:D1)
#   E2 = "Yes"     (bordered, same look as D2/D3)
#   E3 = "Yes"     (bordered, same look as D2/D3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the existing last column (border + header fill)
# into the new column E, the same way a user would copy/paste the column
# formatting before typing the new values in.
$ws.Range("D1:D3").Copy()
$ws.Range("E1:E3").PasteSpecial(-4122)   # xlPasteFormats

# New header cell, matching the style of the other header cells.
$ws.Range("E1").Value = "EXECUTE"

# New data cells.
$ws.Range("E2").Value = "Yes"
$ws.Range("E3").Value = "Yes"

# The data rows (E2:E3) get an explicit (but invisible/white) fill applied
# on top of the border, distinguishing their style from the plain
# border-only style used by D2:D3, while keeping the cells looking
# unfilled/white like the rest of the table.
$ws.Range("E2:E3").Interior.ColorIndex = 2   # white -> visually "no fill"

# Move the active selection, matching the workbook's last saved cursor
# position.
[void]$ws.Range("D6").Select()
